# translations changes before pull from upstream
#
# Insert a new "table_specific_translations" worksheet right after the
# "settings" sheet (and before "model"), and populate it with the
# string_token / text.default / text.spanish translation table.

$wb = $excel.ActiveWorkbook

$settingsSheet = $wb.Worksheets.Item("settings")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $settingsSheet)
$newSheet.Name = "table_specific_translations"

# Header row
$newSheet.Range("A1").Value = "string_token"
$newSheet.Range("B1").Value = "text.default"
$newSheet.Range("C1").Value = "text.spanish"

# Data rows
$newSheet.Range("A2").Value = "address"
$newSheet.Range("B2").Value = "Address"

$newSheet.Range("A3").Value = "id_number"
$newSheet.Range("B3").Value = "ID Number"

$newSheet.Range("A4").Value = "city"
$newSheet.Range("B4").Value = "City"

$newSheet.Range("A5").Value = "telephone"
$newSheet.Range("B5").Value = "Telephon"

$newSheet.Range("A6").Value = "mobile_provider"
$newSheet.Range("B6").Value = "Mobile Provider"

$newSheet.Range("A7").Value = "enable_success"
$newSheet.Range("A8").Value = "disable_success"

$newSheet.Range("B7").Value = "Successfully Enabled!"
$newSheet.Range("B8").Value = "Successfully Disabled!"

$newSheet.Range("A9").Value = "choose_entitlement"
$newSheet.Range("B9").Value = "Choose an Entitlement To Deliver"

$newSheet.Range("A10").Value = "no_entitlements"
$newSheet.Range("B10").Value = "No Entitlements to Deliver"

# Make the newly inserted sheet the active/selected tab, matching the
# upstream workbook's activeTab/tabSelected state.
$newSheet.Activate()
$newSheet.Select()

# Record the last-used cell as the active selection on the new sheet.
$newSheet.Range("B10").Select()
